$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.514.29'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.746.64'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '322.70'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4499'
$ws.Range('E7').Value = '  +5.59%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3501'
$ws.Range('E8').Value = '  -3.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07396'
$ws.Range('E9').Value = '  -1.56%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.70'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.078'
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.43'
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.909'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.096'
$ws.Range('E15').Value = '  -1.77%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.739.66'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '91.19'
$ws.Range('E17').Value = '  -2.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001054'
$ws.Range('E18').Value = '  -1.31%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06395'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.96'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.705'
$ws.Range('E22').Value = '  -3.39%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '27.543.34'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.09'
$ws.Range('E24').Value = '  -1.63%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.119'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '160.83'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.944.20'
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '125.06'
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.025'
$ws.Range('E30').Value = '  -6.33%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.043'
$ws.Range('E31').Value = '  -5.81%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09094'
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.670'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.393'
$ws.Range('E34').Value = '  -3.30%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02266'
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '11.55'
$ws.Range('E36').Value = '  -5.66%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06008'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2061'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.891'
$ws.Range('E39').Value = '  -1.68%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6215'
$ws.Range('E40').Value = '  -2.02%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.176'
$ws.Range('E41').Value = '  -0.88%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.384'
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.693'
$ws.Range('E43').Value = '  -3.56%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.06'
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.694'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5803'
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '121.72'
$ws.Range('E47').Value = '  -1.30%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.924'
$ws.Range('E48').Value = '  -2.48%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.06912'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.114'
$ws.Range('E50').Value = '  -5.53%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '71.18'
$ws.Range('E51').Value = '  -2.87%  '
